$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 49, shifting existing rows 49-69 down to 51-71.
$ws.Range("A49:A50").EntireRow.Insert()

# Populate the two newly inserted rows (49 and 50) with new weekly price data
# for Murcott mandarinas sold in 18kg trays ("bandeja 18 kilos") from Provincia de Limarí.

# Row 49
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).Value = 44468
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100102
$ws.Cells.Item(49, 8).Value = "Cítricos"
$ws.Cells.Item(49, 9).Value = 100102004
$ws.Cells.Item(49, 10).Value = "Mandarina"
$ws.Cells.Item(49, 11).Value = "Murcott"
$ws.Cells.Item(49, 12).Value = "Primera"
$ws.Cells.Item(49, 13).Value = 200
$ws.Cells.Item(49, 14).Value = 7000
$ws.Cells.Item(49, 15).Value = 7500
$ws.Cells.Item(49, 16).Value = 7250
$ws.Cells.Item(49, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(49, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(49, 19).Value = 403
$ws.Cells.Item(49, 20).Value = 18

# Row 50
$ws.Cells.Item(50, 1).Value = 11
$ws.Cells.Item(50, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value = "Bíobío"
$ws.Cells.Item(50, 4).Value = 44468
$ws.Cells.Item(50, 5).Value = 8
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100102
$ws.Cells.Item(50, 8).Value = "Cítricos"
$ws.Cells.Item(50, 9).Value = 100102004
$ws.Cells.Item(50, 10).Value = "Mandarina"
$ws.Cells.Item(50, 11).Value = "Murcott"
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 100
$ws.Cells.Item(50, 14).Value = 6500
$ws.Cells.Item(50, 15).Value = 6500
$ws.Cells.Item(50, 16).Value = 6500
$ws.Cells.Item(50, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(50, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(50, 19).Value = 361
$ws.Cells.Item(50, 20).Value = 18

# Ensure the date cells keep the same date-time number format as the rest of column D.
$ws.Range("D49:D50").NumberFormat = $ws.Range("D51").NumberFormat()

Write-Output "Applied edit"
